$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B31:B40").NumberFormat = "@"
$ws.Range("I31:I40").NumberFormat = "@"
$ws.Cells.Item(31,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(31,2).Value = '000333'
$ws.Cells.Item(31,3).Value = '美的集团'
$ws.Cells.Item(31,4).Value = 2.96
$ws.Cells.Item(31,5).Value = 39.71497483560155
$ws.Cells.Item(31,6).Value = 75.48999999999999
$ws.Cells.Item(31,7).Value = 2998.083450339561
$ws.Cells.Item(31,8).Value = 101173.7496760626
$ws.Cells.Item(31,9).Value = '202506111600'
$ws.Cells.Item(32,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(32,2).Value = '510050'
$ws.Cells.Item(32,3).Value = '上证50ETF'
$ws.Cells.Item(32,4).Value = 4.93
$ws.Cells.Item(32,5).Value = 1808.278443601665
$ws.Cells.Item(32,6).Value = 2.76
$ws.Cells.Item(32,7).Value = 4990.848504340595
$ws.Cells.Item(32,8).Value = 101173.7496760626
$ws.Cells.Item(32,9).Value = '202506111600'
$ws.Cells.Item(33,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(33,2).Value = '510300'
$ws.Cells.Item(33,3).Value = '沪深300ETF'
$ws.Cells.Item(33,4).Value = 4.93
$ws.Cells.Item(33,5).Value = 1247.712126085149
$ws.Cells.Item(33,6).Value = 4
$ws.Cells.Item(33,7).Value = 4990.848504340596
$ws.Cells.Item(33,8).Value = 101173.7496760626
$ws.Cells.Item(33,9).Value = '202506111600'
$ws.Cells.Item(34,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(34,2).Value = '518880'
$ws.Cells.Item(34,3).Value = '黄金ETF'
$ws.Cells.Item(34,4).Value = 4.95
$ws.Cells.Item(34,5).Value = 674.4389870730533
$ws.Cells.Item(34,6).Value = 7.43
$ws.Cells.Item(34,7).Value = 5011.081673952786
$ws.Cells.Item(34,8).Value = 101173.7496760626
$ws.Cells.Item(34,9).Value = '202506111600'
$ws.Cells.Item(35,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(35,2).Value = '600085'
$ws.Cells.Item(35,3).Value = '同仁堂'
$ws.Cells.Item(35,4).Value = 1.97
$ws.Cells.Item(35,5).Value = 52.96734947562319
$ws.Cells.Item(35,6).Value = 37.6
$ws.Cells.Item(35,7).Value = 1991.572340283432
$ws.Cells.Item(35,8).Value = 101173.7496760626
$ws.Cells.Item(35,9).Value = '202506111600'
$ws.Cells.Item(36,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(36,2).Value = '600900'
$ws.Cells.Item(36,3).Value = '长江电力'
$ws.Cells.Item(36,4).Value = 19.86
$ws.Cells.Item(36,5).Value = 665.8903941748626
$ws.Cells.Item(36,6).Value = 30.18
$ws.Cells.Item(36,7).Value = 20096.57209619735
$ws.Cells.Item(36,8).Value = 101173.7496760626
$ws.Cells.Item(36,9).Value = '202506111600'
$ws.Cells.Item(37,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(37,2).Value = '600989'
$ws.Cells.Item(37,3).Value = '宝丰能源'
$ws.Cells.Item(37,4).Value = 4.85
$ws.Cells.Item(37,5).Value = 306.7515982999751
$ws.Cells.Item(37,6).Value = 16
$ws.Cells.Item(37,7).Value = 4908.025572799602
$ws.Cells.Item(37,8).Value = 101173.7496760626
$ws.Cells.Item(37,9).Value = '202506111600'
$ws.Cells.Item(38,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(38,2).Value = 'HK02899'
$ws.Cells.Item(38,3).Value = '紫金矿业'
$ws.Cells.Item(38,4).Value = 20.85
$ws.Cells.Item(38,5).Value = 1106.618293645365
$ws.Cells.Item(38,6).Value = 19.06
$ws.Cells.Item(38,7).Value = 21092.14467688065
$ws.Cells.Item(38,8).Value = 101173.7496760626
$ws.Cells.Item(38,9).Value = '202506111600'
$ws.Cells.Item(39,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(39,2).Value = 'HK06881'
$ws.Cells.Item(39,3).Value = '中国银河'
$ws.Cells.Item(39,4).Value = 5.09
$ws.Cells.Item(39,5).Value = 610.1281790147427
$ws.Cells.Item(39,6).Value = 8.44
$ws.Cells.Item(39,7).Value = 5149.481830884429
$ws.Cells.Item(39,8).Value = 101173.7496760626
$ws.Cells.Item(39,9).Value = '202506111600'
$ws.Cells.Item(40,1).Value = '大智 (稳健智远)'
$ws.Cells.Item(40,2).Value = '100000'
$ws.Cells.Item(40,3).Value = '现金'
$ws.Cells.Item(40,4).Value = 29.6
$ws.Cells.Item(40,5).Value = 29945.09102604357
$ws.Cells.Item(40,6).Value = 1
$ws.Cells.Item(40,7).Value = 29945.09102604357
$ws.Cells.Item(40,8).Value = 101173.7496760626
$ws.Cells.Item(40,9).Value = '202506111600'
$ws.Range("B31:B40").ClearFormats()
$ws.Range("I31:I40").ClearFormats()

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B23:B29").NumberFormat = "@"
$ws.Range("I23:I29").NumberFormat = "@"
$ws.Cells.Item(23,1).Value = '大成 (锐进先锋)'
$ws.Cells.Item(23,2).Value = '000725'
$ws.Cells.Item(23,3).Value = '京东方A'
$ws.Cells.Item(23,4).Value = 5.11
$ws.Cells.Item(23,5).Value = 1243.91705951017
$ws.Cells.Item(23,6).Value = 3.93
$ws.Cells.Item(23,7).Value = 4888.594043874969
$ws.Cells.Item(23,8).Value = 95593.09096603499
$ws.Cells.Item(23,9).Value = '202506111600'
$ws.Cells.Item(24,1).Value = '大成 (锐进先锋)'
$ws.Cells.Item(24,2).Value = '159781'
$ws.Cells.Item(24,3).Value = '科创创业ETF'
$ws.Cells.Item(24,4).Value = 5.06
$ws.Cells.Item(24,5).Value = 9122.058436407913
$ws.Cells.Item(24,6).Value = 0.53
$ws.Cells.Item(24,7).Value = 4834.690971296194
$ws.Cells.Item(24,8).Value = 95593.09096603499
$ws.Cells.Item(24,9).Value = '202506111600'
$ws.Cells.Item(25,1).Value = '大成 (锐进先锋)'
$ws.Cells.Item(25,2).Value = '513100'
$ws.Cells.Item(25,3).Value = '纳指ETF'
$ws.Cells.Item(25,4).Value = 5.19
$ws.Cells.Item(25,5).Value = 3137.523283860047
$ws.Cells.Item(25,6).Value = 1.58
$ws.Cells.Item(25,7).Value = 4957.286788498875
$ws.Cells.Item(25,8).Value = 95593.09096603499
$ws.Cells.Item(25,9).Value = '202506111600'
$ws.Cells.Item(26,1).Value = '大成 (锐进先锋)'
$ws.Cells.Item(26,2).Value = '513290'
$ws.Cells.Item(26,3).Value = '纳指生物科技ETF'
$ws.Cells.Item(26,4).Value = 1.04
$ws.Cells.Item(26,5).Value = 879.627063510763
$ws.Cells.Item(26,6).Value = 1.13
$ws.Cells.Item(26,7).Value = 993.9785817671622
$ws.Cells.Item(26,8).Value = 95593.09096603499
$ws.Cells.Item(26,9).Value = '202506111600'
$ws.Cells.Item(27,1).Value = '大成 (锐进先锋)'
$ws.Cells.Item(27,2).Value = '603119'
$ws.Cells.Item(27,3).Value = '浙江荣泰'
$ws.Cells.Item(27,4).Value = 43.67
$ws.Cells.Item(27,5).Value = 1051.546584462582
$ws.Cells.Item(27,6).Value = 39.7
$ws.Cells.Item(27,7).Value = 41746.39940316451
$ws.Cells.Item(27,8).Value = 95593.09096603499
$ws.Cells.Item(27,9).Value = '202506111600'
$ws.Cells.Item(28,1).Value = '大成 (锐进先锋)'
$ws.Cells.Item(28,2).Value = '688290'
$ws.Cells.Item(28,3).Value = '景业智能'
$ws.Cells.Item(28,4).Value = 9.01
$ws.Cells.Item(28,5).Value = 161.2702946560293
$ws.Cells.Item(28,6).Value = 53.43
$ws.Cells.Item(28,7).Value = 8616.671843471646
$ws.Cells.Item(28,8).Value = 95593.09096603499
$ws.Cells.Item(28,9).Value = '202506111600'
$ws.Cells.Item(29,1).Value = '大成 (锐进先锋)'
$ws.Cells.Item(29,2).Value = '100000'
$ws.Cells.Item(29,3).Value = '现金'
$ws.Cells.Item(29,4).Value = 30.92
$ws.Cells.Item(29,5).Value = 29555.46933396164
$ws.Cells.Item(29,6).Value = 1
$ws.Cells.Item(29,7).Value = 29555.46933396164
$ws.Cells.Item(29,8).Value = 95593.09096603499
$ws.Cells.Item(29,9).Value = '202506111600'
$ws.Range("B23:B29").ClearFormats()
$ws.Range("I23:I29").ClearFormats()

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("B46:B60").NumberFormat = "@"
$ws.Range("I46:I60").NumberFormat = "@"
$ws.Cells.Item(46,1).Value = '范式进化投资组合'
$ws.Cells.Item(46,2).Value = '000333'
$ws.Cells.Item(46,3).Value = '美的集团'
$ws.Cells.Item(46,4).Value = 1
$ws.Cells.Item(46,5).Value = 13.2658076636599
$ws.Cells.Item(46,6).Value = 75.48999999999999
$ws.Cells.Item(46,7).Value = 1001.435820529686
$ws.Cells.Item(46,8).Value = 99865.51084235184
$ws.Cells.Item(46,9).Value = '202506111600'
$ws.Cells.Item(47,1).Value = '范式进化投资组合'
$ws.Cells.Item(47,2).Value = '000725'
$ws.Cells.Item(47,3).Value = '京东方A'
$ws.Cells.Item(47,4).Value = 4.97
$ws.Cells.Item(47,5).Value = 1262.931689191865
$ws.Cells.Item(47,6).Value = 3.93
$ws.Cells.Item(47,7).Value = 4963.321538524029
$ws.Cells.Item(47,8).Value = 99865.51084235184
$ws.Cells.Item(47,9).Value = '202506111600'
$ws.Cells.Item(48,1).Value = '范式进化投资组合'
$ws.Cells.Item(48,2).Value = '159781'
$ws.Cells.Item(48,3).Value = '科创创业ETF'
$ws.Cells.Item(48,4).Value = 4.92
$ws.Cells.Item(48,5).Value = 9261.499054073673
$ws.Cells.Item(48,6).Value = 0.53
$ws.Cells.Item(48,7).Value = 4908.594498659047
$ws.Cells.Item(48,8).Value = 99865.51084235184
$ws.Cells.Item(48,9).Value = '202506111600'
$ws.Cells.Item(49,1).Value = '范式进化投资组合'
$ws.Cells.Item(49,2).Value = '510050'
$ws.Cells.Item(49,3).Value = '上证50ETF'
$ws.Cells.Item(49,4).Value = 5.01
$ws.Cells.Item(49,5).Value = 1812.03242362311
$ws.Cells.Item(49,6).Value = 2.76
$ws.Cells.Item(49,7).Value = 5001.209489199783
$ws.Cells.Item(49,8).Value = 99865.51084235184
$ws.Cells.Item(49,9).Value = '202506111600'
$ws.Cells.Item(50,1).Value = '范式进化投资组合'
$ws.Cells.Item(50,2).Value = '510300'
$ws.Cells.Item(50,3).Value = '沪深300ETF'
$ws.Cells.Item(50,4).Value = 5.01
$ws.Cells.Item(50,5).Value = 1250.302372299946
$ws.Cells.Item(50,6).Value = 4
$ws.Cells.Item(50,7).Value = 5001.209489199784
$ws.Cells.Item(50,8).Value = 99865.51084235184
$ws.Cells.Item(50,9).Value = '202506111600'
$ws.Cells.Item(51,1).Value = '范式进化投资组合'
$ws.Cells.Item(51,2).Value = '513100'
$ws.Cells.Item(51,3).Value = '纳指ETF'
$ws.Cells.Item(51,4).Value = 1.01
$ws.Cells.Item(51,5).Value = 637.096750216533
$ws.Cells.Item(51,6).Value = 1.58
$ws.Cells.Item(51,7).Value = 1006.612865342122
$ws.Cells.Item(51,8).Value = 99865.51084235184
$ws.Cells.Item(51,9).Value = '202506111600'
$ws.Cells.Item(52,1).Value = '范式进化投资组合'
$ws.Cells.Item(52,2).Value = '513290'
$ws.Cells.Item(52,3).Value = '纳指生物科技ETF'
$ws.Cells.Item(52,4).Value = 1.01
$ws.Cells.Item(52,5).Value = 893.0731230713899
$ws.Cells.Item(52,6).Value = 1.13
$ws.Cells.Item(52,7).Value = 1009.17262907067
$ws.Cells.Item(52,8).Value = 99865.51084235184
$ws.Cells.Item(52,9).Value = '202506111600'
$ws.Cells.Item(53,1).Value = '范式进化投资组合'
$ws.Cells.Item(53,2).Value = '518880'
$ws.Cells.Item(53,3).Value = '黄金ETF'
$ws.Cells.Item(53,4).Value = 1.01
$ws.Cells.Item(53,5).Value = 135.1678240324266
$ws.Cells.Item(53,6).Value = 7.43
$ws.Cells.Item(53,7).Value = 1004.29693256093
$ws.Cells.Item(53,8).Value = 99865.51084235184
$ws.Cells.Item(53,9).Value = '202506111600'
$ws.Cells.Item(54,1).Value = '范式进化投资组合'
$ws.Cells.Item(54,2).Value = '600085'
$ws.Cells.Item(54,3).Value = '同仁堂'
$ws.Cells.Item(54,4).Value = 1
$ws.Cells.Item(54,5).Value = 26.53865475829018
$ws.Cells.Item(54,6).Value = 37.6
$ws.Cells.Item(54,7).Value = 997.8534189117107
$ws.Cells.Item(54,8).Value = 99865.51084235184
$ws.Cells.Item(54,9).Value = '202506111600'
$ws.Cells.Item(55,1).Value = '范式进化投资组合'
$ws.Cells.Item(55,2).Value = '600900'
$ws.Cells.Item(55,3).Value = '长江电力'
$ws.Cells.Item(55,4).Value = 1.01
$ws.Cells.Item(55,5).Value = 33.36363902067901
$ws.Cells.Item(55,6).Value = 30.18
$ws.Cells.Item(55,7).Value = 1006.914625644093
$ws.Cells.Item(55,8).Value = 99865.51084235184
$ws.Cells.Item(55,9).Value = '202506111600'
$ws.Cells.Item(56,1).Value = '范式进化投资组合'
$ws.Cells.Item(56,2).Value = '600989'
$ws.Cells.Item(56,3).Value = '宝丰能源'
$ws.Cells.Item(56,4).Value = 4.92
$ws.Cells.Item(56,5).Value = 307.3884135955614
$ws.Cells.Item(56,6).Value = 16
$ws.Cells.Item(56,7).Value = 4918.214617528983
$ws.Cells.Item(56,8).Value = 99865.51084235184
$ws.Cells.Item(56,9).Value = '202506111600'
$ws.Cells.Item(57,1).Value = '范式进化投资组合'
$ws.Cells.Item(57,2).Value = '603119'
$ws.Cells.Item(57,3).Value = '浙江荣泰'
$ws.Cells.Item(57,4).Value = 0.9399999999999999
$ws.Cells.Item(57,5).Value = 23.72490270018873
$ws.Cells.Item(57,6).Value = 39.7
$ws.Cells.Item(57,7).Value = 941.8786371974926
$ws.Cells.Item(57,8).Value = 99865.51084235184
$ws.Cells.Item(57,9).Value = '202506111600'
$ws.Cells.Item(58,1).Value = '范式进化投资组合'
$ws.Cells.Item(58,2).Value = 'HK02899'
$ws.Cells.Item(58,3).Value = '紫金矿业'
$ws.Cells.Item(58,4).Value = 1.06
$ws.Cells.Item(58,5).Value = 55.44578147671601
$ws.Cells.Item(58,6).Value = 19.06
$ws.Cells.Item(58,7).Value = 1056.796594946207
$ws.Cells.Item(58,8).Value = 99865.51084235184
$ws.Cells.Item(58,9).Value = '202506111600'
$ws.Cells.Item(59,1).Value = '范式进化投资组合'
$ws.Cells.Item(59,2).Value = 'HK06881'
$ws.Cells.Item(59,3).Value = '中国银河'
$ws.Cells.Item(59,4).Value = 1.03
$ws.Cells.Item(59,5).Value = 122.2789606161316
$ws.Cells.Item(59,6).Value = 8.44
$ws.Cells.Item(59,7).Value = 1032.034427600151
$ws.Cells.Item(59,8).Value = 99865.51084235184
$ws.Cells.Item(59,9).Value = '202506111600'
$ws.Cells.Item(60,1).Value = '范式进化投资组合'
$ws.Cells.Item(60,2).Value = '100000'
$ws.Cells.Item(60,3).Value = '现金'
$ws.Cells.Item(60,4).Value = 66.09999999999999
$ws.Cells.Item(60,5).Value = 66015.96525743716
$ws.Cells.Item(60,6).Value = 1
$ws.Cells.Item(60,7).Value = 66015.96525743716
$ws.Cells.Item(60,8).Value = 99865.51084235184
$ws.Cells.Item(60,9).Value = '202506111600'
$ws.Range("B46:B60").ClearFormats()
$ws.Range("I46:I60").ClearFormats()

